$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "BENCHMARKING-2026"
$ws.Range("D2").Value = "23,81 TL - 23,81 TL"
$ws.Range("C3").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("E3").Value = ""
$ws.Range("I3").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("K3").Value = ""
$ws.Range("C4").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("E4").Value = ""
$ws.Range("I4").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("K4").Value = ""
$ws.Range("C5").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("E5").Value = ""
$ws.Range("I5").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("K5").Value = ""
$ws.Range("C6").Value = "6,09 TL - 12,19 TL - 152,35 TL"
$ws.Range("E6").Value = ""
$ws.Range("I6").Value = "6,09 TL - 12,19 TL - 152,35 TL"
$ws.Range("K6").Value = ""
$ws.Range("D7").Value = "%1,6"
$ws.Range("C8").Value = "14,29 TL - 28,57 TL - 300 TL"
$ws.Range("E8").Value = ""
$ws.Range("I8").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("K8").Value = ""
$ws.Range("C9").Value = "14,29 TL - 28,57 TL - 300 TL"
$ws.Range("E9").Value = ""
$ws.Range("I9").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("K9").Value = ""
$ws.Range("C10").Value = "14,29 TL - 28,57 TL - 300 TL"
$ws.Range("E10").Value = ""
$ws.Range("I10").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("K10").Value = ""
$ws.Range("C11").Value = "3,04 TL - 6,09 TL - 76,17 TL"
$ws.Range("E11").Value = ""
$ws.Range("I11").Value = "3,04 TL - 6,09 TL - 76,17 TL"
$ws.Range("K11").Value = ""
$ws.Range("C12").Value = "WU: 1.000,01 USD–9,51 USD"
$ws.Range("K12").Value = ""
$ws.Range("C13").Value = "Hesaba: Asgari 0 TL | Azami 0,94 TL"
$ws.Range("E13").Value = ""
$ws.Range("I13").Value = "Hesaba: Asgari 1 TL | Azami 6,09 TL"
$ws.Range("K13").Value = ""
$ws.Range("C14").Value = "40.000 TL - 1.904,76 TL"
$ws.Range("E14").Value = ""
$ws.Range("K14").Value = ""
$ws.Range("I23").Value = "350 TL"
